# Deploy the implementation guide.
# Update the build Date on the Metadata sheet and fix the casing of the
# "Missing - ..." Display values on the Concepts sheet.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# Update the generation date/time.
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# Fix capitalization of the Display column (column C) for the "Missing - ..."
# concepts, while leaving the Code column (column B) untouched.
$concepts.Range("C3").Value = "Missing - Restricted Access"
$concepts.Range("C4").Value = "Missing - Not Provided"
$concepts.Range("C5").Value = "Missing - Not Collected"
